# Governance Body Suite - test-run bookkeeping update
#
# - "Test Cases" sheet: row 5 (the GBUpdate test case) Runmode flips from
#   "Y" to "N", and its stale "FAIL" Result is cleared out.
# - "Test Cases" becomes the active/selected sheet (was "GBCreation"),
#   with the selection resting on C4.

$wb = $excel.ActiveWorkbook

$wsTestCases = $wb.Worksheets.Item("Test Cases")

# Flip the Runmode for the GBUpdate row to "N" and clear the old Result.
$wsTestCases.Range("C5").Value = "N"
$wsTestCases.Range("D5").Value = $null

# Make "Test Cases" the active sheet/tab, selecting C4 (this also drops
# the previous tab-selected state from "GBCreation").
[void]$wsTestCases.Activate()
[void]$wsTestCases.Range("C4").Select()
